$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Determine the last used row in the sheet (data starts at row 2)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 530 }

# Column C holds the "Förändrad" (Changed) date. Update every data row
# from the old date serial 45179 (2023-09-10) to the new date serial
# 45180 (2023-09-11).
$range = $ws.Range("C2:C$lastRow")
$range.Value2 = 45180
